$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meloxicam")

# D2 previously held the text "C0083381"; replace it with the numeric value 41493
$ws.Range("D2").Value = 41493

# Update the active selection to D13, matching the saved workbook state
$ws.Range("D13").Select()
